$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.066677059972287
$ws.Range("D2").Value = 1.049718084173496
$ws.Range("E2").Value = 1.070899531701963
$ws.Range("F2").Value = 1.079845335107357
$ws.Range("I2").Value = 1.048797781645715
$ws.Range("J2").Value = 1.071625841969085
$ws.Range("K2").Value = 1.052474363342899
$ws.Range("L2").Value = 1.073598139330804
$ws.Range("M2").Value = 1.082520319984118
$ws.Range("N2").Value = 1.073147673581883

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.068159856152864
$ws.Range("D3").Value = 1.050443920519322
$ws.Range("E3").Value = 1.072224556851048
$ws.Range("F3").Value = 1.081260699795537
$ws.Range("I3").Value = 1.049175282993491
$ws.Range("J3").Value = 1.072762203358879
$ws.Range("K3").Value = 1.053012725120281
$ws.Range("L3").Value = 1.074738357812207
$ws.Range("M3").Value = 1.083752380917148
$ws.Range("N3").Value = 1.074285648735193

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.069118236624402
$ws.Range("D4").Value = 1.050913074677537
$ws.Range("E4").Value = 1.073081168872642
$ws.Range("F4").Value = 1.082175901594509
$ws.Range("I4").Value = 1.049417907112838
$ws.Range("J4").Value = 1.073495973158782
$ws.Range("K4").Value = 1.053359901665485
$ws.Range("L4").Value = 1.075474839556287
$ws.Range("M4").Value = 1.084548439828315
$ws.Range("N4").Value = 1.075020460572377

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.069520885688234
$ws.Range("D5").Value = 1.051110185739133
$ws.Range("E5").Value = 1.073441109568662
$ws.Range("F5").Value = 1.082560505838401
$ws.Range("I5").Value = 1.049519514024407
$ws.Range("J5").Value = 1.073804087836347
$ws.Range("K5").Value = 1.053505573461551
$ws.Range("L5").Value = 1.075784145908344
$ws.Range("M5").Value = 1.084882828325896
$ws.Range("N5").Value = 1.075329012808125

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.069588477521726
$ws.Range("D6").Value = 1.051143274462792
$ws.Range("E6").Value = 1.073501534838695
$ws.Range("F6").Value = 1.082625074135989
$ws.Range("I6").Value = 1.049536551310471
$ws.Range("J6").Value = 1.07385580054387
$ws.Range("K6").Value = 1.053530015939412
$ws.Range("L6").Value = 1.075836061722517
$ws.Range("M6").Value = 1.084938957655181
$ws.Range("N6").Value = 1.075380798953624

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.069123617831426
$ws.Range("D7").Value = 1.050915708960947
$ws.Range("E7").Value = 1.073085979115194
$ws.Range("F7").Value = 1.082181041264873
$ws.Range("I7").Value = 1.049419266328946
$ws.Range("J7").Value = 1.073500091619455
$ws.Range("K7").Value = 1.053361849242751
$ws.Range("L7").Value = 1.075478973737739
$ws.Range("M7").Value = 1.084552909017237
$ws.Range("N7").Value = 1.075024584881736

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.06717840694397
$ws.Range("D8").Value = 1.049963490192141
$ws.Range("E8").Value = 1.071347492055831
$ws.Range("F8").Value = 1.080323798562789
$ws.Range("I8").Value = 1.048925701359133
$ws.Range("J8").Value = 1.072010200619582
$ws.Range("K8").Value = 1.052656550154495
$ws.Range("L8").Value = 1.073983756644204
$ws.Range("M8").Value = 1.082936945187658
$ws.Range("N8").Value = 1.07353257806575

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.063742081039119
$ws.Range("D9").Value = 1.048281609433723
$ws.Range("E9").Value = 1.068277957970566
$ws.Range("F9").Value = 1.077046014570188
$ws.Range("I9").Value = 1.048043317572139
$ws.Range("J9").Value = 1.069372881415453
$ws.Range("K9").Value = 1.051404636734577
$ws.Range("L9").Value = 1.071338715378013
$ws.Range("M9").Value = 1.080080262680696
$ws.Range("N9").Value = 1.070891513565848

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.061445004874942
$ws.Range("D10").Value = 1.047157641360188
$ws.Range("E10").Value = 1.066227196693296
$ws.Range("F10").Value = 1.074857081327321
$ws.Range("I10").Value = 1.047446458993285
$ws.Range("J10").Value = 1.067606348065932
$ws.Range("K10").Value = 1.050563842998931
$ws.Range("L10").Value = 1.069568163056478
$ws.Range("M10").Value = 1.078169360550319
$ws.Range("N10").Value = 1.069122471536345

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.060448787423283
$ws.Range("D11").Value = 1.046670293008008
$ws.Range("E11").Value = 1.06533808096182
$ws.Range("F11").Value = 1.073908288100458
$ws.Range("I11").Value = 1.047185951435095
$ws.Range("J11").Value = 1.066839380500712
$ws.Range("K11").Value = 1.050198285571038
$ws.Range("L11").Value = 1.068799725437121
$ws.Range("M11").Value = 1.077340323022421
$ws.Range("N11").Value = 1.068354414789211

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.060078504504884
$ws.Range("D12").Value = 1.04648916925521
$ws.Range("E12").Value = 1.065007649145361
$ws.Range("F12").Value = 1.073555711847123
$ws.Range("I12").Value = 1.047088875405748
$ws.Range("J12").Value = 1.066554181326725
$ws.Range("K12").Value = 1.050062276091064
$ws.Range("L12").Value = 1.068514021021234
$ws.Range("M12").Value = 1.077032134717863
$ws.Range("N12").Value = 1.068068810599684

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.060157942587632
$ws.Range("D13").Value = 1.046528025501014
$ws.Range("E13").Value = 1.065078535891423
$ws.Range("F13").Value = 1.073631347640571
$ws.Range("I13").Value = 1.047109712688647
$ws.Range("J13").Value = 1.066615371786378
$ws.Range("K13").Value = 1.050091460801579
$ws.Range("L13").Value = 1.068575318007258
$ws.Range("M13").Value = 1.077098253419983
$ws.Range("N13").Value = 1.068130087956808

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.060418184732199
$ws.Range("D14").Value = 1.046655323327549
$ws.Range("E14").Value = 1.065310770962327
$ws.Range("F14").Value = 1.073879147179274
$ws.Range("I14").Value = 1.047177933474128
$ws.Range("J14").Value = 1.06681581227198
$ws.Range("K14").Value = 1.050187047587758
$ws.Range("L14").Value = 1.06877611461138
$ws.Range("M14").Value = 1.077314853150611
$ws.Range("N14").Value = 1.06833081309089

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.060578496020105
$ws.Range("D15").Value = 1.046733742313478
$ws.Range("E15").Value = 1.065453835331982
$ws.Range("F15").Value = 1.074031804293747
$ws.Range("I15").Value = 1.047219925163137
$ws.Range("J15").Value = 1.066939268586792
$ws.Range("K15").Value = 1.050245911875744
$ws.Range("L15").Value = 1.068899795755457
$ws.Range("M15").Value = 1.077448274507852
$ws.Range("N15").Value = 1.06845444472783

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.061511086840823
$ws.Range("D16").Value = 1.047189970990926
$ws.Range("E16").Value = 1.066286180214372
$ws.Range("F16").Value = 1.074920028615073
$ws.Range("I16").Value = 1.047463704381135
$ws.Range("J16").Value = 1.067657205561199
$ws.Range("K16").Value = 1.050588072351615
$ws.Range("L16").Value = 1.069619123824751
$ws.Range("M16").Value = 1.078224346743577
$ws.Range("N16").Value = 1.069173401255089

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.062095651064885
$ws.Range("D17").Value = 1.047475972833434
$ws.Range("E17").Value = 1.066807983407327
$ws.Range("F17").Value = 1.075476923900929
$ws.Range("I17").Value = 1.047616066788058
$ws.Range("J17").Value = 1.06810699638742
$ws.Range("K17").Value = 1.050802301029344
$ws.Range("L17").Value = 1.070069859732073
$ws.Range("M17").Value = 1.07871072273638
$ws.Range("N17").Value = 1.069623830835875

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.062436466695764
$ws.Range("D18").Value = 1.047642729012844
$ws.Range("E18").Value = 1.067112234477951
$ws.Range("F18").Value = 1.075801658275116
$ws.Range("I18").Value = 1.04770473819646
$ws.Range("J18").Value = 1.068369154653997
$ws.Range("K18").Value = 1.050927113512715
$ws.Range("L18").Value = 1.070332595303452
$ws.Range("M18").Value = 1.078994263229217
$ws.Range("N18").Value = 1.069886361397261

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.062552650675948
$ws.Range("D19").Value = 1.047699577768235
$ws.Range("E19").Value = 1.06721595810852
$ws.Range("F19").Value = 1.075912368684088
$ws.Range("I19").Value = 1.047734939156062
$ws.Range("J19").Value = 1.068458510587116
$ws.Range("K19").Value = 1.050969647055017
$ws.Range("L19").Value = 1.070422152495656
$ws.Range("M19").Value = 1.079090917197232
$ws.Range("N19").Value = 1.069975844226054

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.062032948509377
$ws.Range("D20").Value = 1.047445294156962
$ws.Range("E20").Value = 1.066752010067529
$ws.Range("F20").Value = 1.075417183996031
$ws.Range("I20").Value = 1.047599740341317
$ws.Range("J20").Value = 1.068058758518476
$ws.Range("K20").Value = 1.050779331172068
$ws.Range("L20").Value = 1.070021517768371
$ws.Range("M20").Value = 1.07865855520748
$ws.Range("N20").Value = 1.069575524463624

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.060341556688742
$ws.Range("D21").Value = 1.046617840085683
$ws.Range("E21").Value = 1.065242388391666
$ws.Range("F21").Value = 1.073806180644298
$ws.Range("I21").Value = 1.047157852778714
$ws.Range("J21").Value = 1.066756796238897
$ws.Range("K21").Value = 1.050158905890457
$ws.Range("L21").Value = 1.068716992582259
$ws.Range("M21").Value = 1.077251076789906
$ws.Range("N21").Value = 1.068271713248271

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.059276699392639
$ws.Range("D22").Value = 1.046097002573034
$ws.Range("E22").Value = 1.064292215882164
$ws.Range("F22").Value = 1.072792395106592
$ws.Range("I22").Value = 1.046878214910738
$ws.Range("J22").Value = 1.065936385519882
$ws.Range("K22").Value = 1.049767516298321
$ws.Range("L22").Value = 1.067895206632431
$ws.Range("M22").Value = 1.076364708111923
$ws.Range("N22").Value = 1.067450137451932

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.059841336855347
$ws.Range("D23").Value = 1.046373164217868
$ws.Range("E23").Value = 1.064796018380807
$ws.Range("F23").Value = 1.073329907911199
$ws.Range("I23").Value = 1.047026628009072
$ws.Range("J23").Value = 1.066371474764488
$ws.Range("K23").Value = 1.049975123427337
$ws.Range("L23").Value = 1.068331002456592
$ws.Range("M23").Value = 1.0768347265284
$ws.Range("N23").Value = 1.067885844573171

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.062061281561453
$ws.Range("D24").Value = 1.047459156727965
$ws.Range("E24").Value = 1.066777302344193
$ws.Range("F24").Value = 1.075444178175992
$ws.Range("I24").Value = 1.047607118174899
$ws.Range("J24").Value = 1.068080555744425
$ws.Range("K24").Value = 1.050789710705429
$ws.Range("L24").Value = 1.07004336194976
$ws.Range("M24").Value = 1.07868212794565
$ws.Range("N24").Value = 1.069597352644135

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.064631515655299
$ws.Range("D25").Value = 1.048716889950914
$ws.Range("E25").Value = 1.069072260051892
$ws.Range("F25").Value = 1.077894036244381
$ws.Range("I25").Value = 1.048272944479337
$ws.Range("J25").Value = 1.070056136670988
$ws.Range("K25").Value = 1.051729370869115
$ws.Range("L25").Value = 1.072023768762897
$ws.Range("M25").Value = 1.080819898071733
$ws.Range("N25").Value = 1.071575739122215
